$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table of (cell address -> new text value) derived from the commit diff.
# Values must remain TEXT (they mirror the source site's inline strings, e.g. "321.16", "6.16%"),
# so each cell is pre-formatted as Text ("@") before the value is written; this stops Excel
# from reinterpreting a numeric-looking string (or a "N.NN%" string) as a Number/Percentage.
$changes = @(
    @{Addr="D2"; Value="321.16"},
    @{Addr="E2"; Value="6.16%"},
    @{Addr="D3"; Value="49.13"},
    @{Addr="E3"; Value="11.24%"},
    @{Addr="D4"; Value="5.325"},
    @{Addr="E4"; Value="4.72%"},
    @{Addr="D5"; Value="0.08079"},
    @{Addr="E5"; Value="4.86%"},
    @{Addr="D6"; Value="4.596"},
    @{Addr="E6"; Value="4.09%"},
    @{Addr="D7"; Value="1.340"},
    @{Addr="E7"; Value="28.07%"},
    @{Addr="D8"; Value="1.639"},
    @{Addr="E8"; Value="1.44%"},
    @{Addr="D9"; Value="0.1278"},
    @{Addr="E9"; Value="0.14%"},
    @{Addr="D10"; Value="0.1967"},
    @{Addr="E10"; Value="5.70%"},
    @{Addr="D11"; Value="0.09669"},
    @{Addr="E11"; Value="4.65%"},
    @{Addr="D12"; Value="0.04716"},
    @{Addr="E12"; Value="12.59%"},
    @{Addr="D13"; Value="0.1046"},
    @{Addr="E13"; Value="0.02%"},
    @{Addr="D14"; Value="0.001327"},
    @{Addr="E14"; Value="3.59%"},
    @{Addr="D15"; Value="0.04196"},
    @{Addr="E15"; Value="0.14%"},
    @{Addr="D16"; Value="0.005783"},
    @{Addr="E16"; Value="0.52%"},
    @{Addr="D17"; Value="3.341"},
    @{Addr="E17"; Value="-0.17%"},
    @{Addr="D18"; Value="2.442"},
    @{Addr="E18"; Value="4.78%"},
    @{Addr="D19"; Value="0.3507"},
    @{Addr="E19"; Value="4.92%"},
    @{Addr="D20"; Value="8.013"},
    @{Addr="E20"; Value="-0.94%"},
    @{Addr="E21"; Value="0.30%"},
    @{Addr="D23"; Value="0.001311"},
    @{Addr="E23"; Value="2.22%"},
    @{Addr="D24"; Value="0.004334"},
    @{Addr="E24"; Value="-1.82%"},
    @{Addr="D25"; Value="0.0001349"},
    @{Addr="E25"; Value="0.02%"},
    @{Addr="D26"; Value="0.0003533"},
    @{Addr="E26"; Value="-95.28%"},
    @{Addr="D38"; Value="0.02727"},
    @{Addr="E38"; Value="9.43%"},
    @{Addr="D39"; Value="0.06084"},
    @{Addr="E39"; Value="14.76%"},
    @{Addr="D40"; Value="0.01084"},
    @{Addr="E40"; Value="82.99%"},
    @{Addr="D41"; Value="0.008025"},
    @{Addr="E41"; Value="3.91%"},
    @{Addr="E42"; Value="8.50%"},
    @{Addr="D43"; Value="0.007884"},
    @{Addr="E43"; Value="7.09%"},
    @{Addr="D44"; Value="0.008651"},
    @{Addr="E44"; Value="14.79%"},
    @{Addr="D45"; Value="0.3501"},
    @{Addr="E45"; Value="16.20%"},
    @{Addr="D46"; Value="0.00006843"},
    @{Addr="E46"; Value="2.94%"},
    @{Addr="D47"; Value="0.00000000749"},
    @{Addr="E47"; Value="-0.09%"},
    @{Addr="D48"; Value="0.05915"},
    @{Addr="E48"; Value="37.28%"},
    @{Addr="D49"; Value="0.003993"},
    @{Addr="E49"; Value="-4.94%"},
    @{Addr="D50"; Value="0.00002096"},
    @{Addr="E50"; Value="-0.09%"},
    @{Addr="D51"; Value="0.0001996"},
    @{Addr="E51"; Value="-0.09%"}
)

foreach ($chg in $changes) {
    $rng = $ws.Range($chg.Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $chg.Value
}

